$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Mapping of (row, col) -> new value. Row/Col are 1-indexed Word table
# coordinates. Only the five "content" rows (1, 5, 9, 13, 17) hold values;
# the other rows are blank spacer rows.
$updates = @(
    @{ Row = 1;  Col = 1; Text = "44÷4=" },
    @{ Row = 1;  Col = 2; Text = "19÷5=" },
    @{ Row = 1;  Col = 3; Text = "77÷8=" },
    @{ Row = 1;  Col = 4; Text = "28÷6=" },
    @{ Row = 1;  Col = 5; Text = "68÷9=" },

    @{ Row = 5;  Col = 1; Text = "68÷2=" },
    @{ Row = 5;  Col = 2; Text = "84÷5=" },
    @{ Row = 5;  Col = 3; Text = "62÷9=" },
    @{ Row = 5;  Col = 4; Text = "84÷2=" },
    @{ Row = 5;  Col = 5; Text = "75÷4=" },

    @{ Row = 9;  Col = 1; Text = "78÷8=" },
    @{ Row = 9;  Col = 2; Text = "24÷9=" },
    @{ Row = 9;  Col = 3; Text = "27÷3=" },
    @{ Row = 9;  Col = 4; Text = "48÷3=" },
    @{ Row = 9;  Col = 5; Text = "66÷5=" },

    @{ Row = 13; Col = 1; Text = "22÷3=" },
    @{ Row = 13; Col = 2; Text = "47÷8=" },
    @{ Row = 13; Col = 3; Text = "73÷4=" },
    @{ Row = 13; Col = 4; Text = "26÷5=" },
    @{ Row = 13; Col = 5; Text = "77÷7=" },

    @{ Row = 17; Col = 1; Text = "79÷5=" },
    @{ Row = 17; Col = 2; Text = "59÷2=" },
    @{ Row = 17; Col = 3; Text = "74÷7=" },
    @{ Row = 17; Col = 4; Text = "65÷9=" },
    @{ Row = 17; Col = 5; Text = "76÷2=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $u.Text
}
